# Refresh cached market-board price / profit figures on the Garuda_Profits sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107: Another Man's Ink
$ws.Range("H107").Value = 146.58824
$ws.Range("I107").Value = 150.13333
$ws.Range("J107").Value = 120
$ws.Range("K107").Value = 150.13333
$ws.Range("L107").Value = 120
$ws.Range("M107").Value = 1769.86667
$ws.Range("N107").Value = -3960

# Row 129: Practical Command
$ws.Range("H129").Value = 28118.594
$ws.Range("I129").Value = 528.7646999999999
$ws.Range("J129").Value = 51569.95
$ws.Range("K129").Value = 1586.2941
$ws.Range("L129").Value = 154709.85
$ws.Range("M129").Value = 3413.7059
$ws.Range("N129").Value = -164709.85

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth
$ws.Range("H5").Value = 55556390
$ws.Range("I5").Value = 83333840
$ws.Range("K5").Value = 83333840
$ws.Range("M5").Value = -83333728

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 1030.0416
$ws.Range("I45").Value = 936.3125
$ws.Range("J45").Value = 1217.5
$ws.Range("K45").Value = 936.3125
$ws.Range("L45").Value = 1217.5
$ws.Range("M45").Value = -559.3125
$ws.Range("N45").Value = -1971.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences
$ws.Range("H4").Value = 55556390
$ws.Range("I4").Value = 83333840
$ws.Range("K4").Value = 83333840
$ws.Range("M4").Value = -83333725

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 3211.889
$ws.Range("I105").Value = 3215.2856
$ws.Range("K105").Value = 3215.2856
$ws.Range("M105").Value = -1468.2856

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1168.04
$ws.Range("I58").Value = 879.3158
$ws.Range("J58").Value = 2082.3333
$ws.Range("K58").Value = 879.3158
$ws.Range("L58").Value = 2082.3333
$ws.Range("M58").Value = -676.3158
$ws.Range("N58").Value = -2488.3333

# Row 136: Turali Quality
$ws.Range("H136").Value = 1168.04
$ws.Range("I136").Value = 879.3158
$ws.Range("J136").Value = 2082.3333
$ws.Range("K136").Value = 2637.9474
$ws.Range("L136").Value = 6246.999899999999
$ws.Range("M136").Value = -87.94740000000002
$ws.Range("N136").Value = -11346.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 25: Flakes for Friends
$ws.Range("H25").Value = 500
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1500
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value = 0
$ws.Range("M25").Value = -1331

# Row 30: Picnic Panic
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1500
$ws.Range("L30").ClearContents()
$ws.Range("N30").Value = 0
$ws.Range("M30").Value = -1398

# Row 38: Pretty as a Picture
$ws.Range("H38").Value = 296
$ws.Range("I38").Value = 46.666668
$ws.Range("J38").Value = 445.6
$ws.Range("K38").Value = 140.000004
$ws.Range("L38").Value = 1336.8
$ws.Range("M38").Value = 206.999996
$ws.Range("N38").Value = -2030.8

# Row 104: Fits to a Tea
$ws.Range("H104").Value = 2985
$ws.Range("I104").Value = 1963
$ws.Range("J104").Value = 5029
$ws.Range("K104").Value = 5889
$ws.Range("L104").Value = 15087
$ws.Range("M104").Value = -3268
$ws.Range("N104").Value = -20329

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 2774.0942
$ws.Range("I131").Value = 4972.7085
$ws.Range("J131").Value = 954.5517
$ws.Range("K131").Value = 14918.1255
$ws.Range("L131").Value = 2863.6551
$ws.Range("M131").Value = -9878.125499999998
$ws.Range("N131").Value = -12943.6551

$ws = $wb.Worksheets.Item("GSM")
# Row 46: Burning the Midnight Oil
$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -9844

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 3887.3635
$ws.Range("I80").Value = 3013.75
$ws.Range("J80").Value = 4386.5713
$ws.Range("K80").Value = 3013.75
$ws.Range("L80").Value = 4386.5713
$ws.Range("M80").Value = -2015.75
$ws.Range("N80").Value = -6382.5713

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 3887.3635
$ws.Range("I83").Value = 3013.75
$ws.Range("J83").Value = 4386.5713
$ws.Range("K83").Value = 15068.75
$ws.Range("L83").Value = 21932.8565
$ws.Range("M83").Value = -10076.75
$ws.Range("N83").Value = -31916.8565

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 62504196
$ws.Range("I122").Value = 200006830
$ws.Range("J122").Value = 2996.182
$ws.Range("K122").Value = 600020490
$ws.Range("L122").Value = 8988.545999999998
$ws.Range("M122").Value = -600018040
$ws.Range("N122").Value = -13888.546

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 2282.0527
$ws.Range("I126").Value = 2417.4666
$ws.Range("K126").Value = 7252.399800000001
$ws.Range("M126").Value = -4782.399800000001

# Row 132: On Board for Lar
$ws.Range("H132").Value = 2297.2058
$ws.Range("I132").Value = 2252.9167
$ws.Range("J132").Value = 2403.5
$ws.Range("K132").Value = 6758.750100000001
$ws.Range("L132").Value = 7210.5
$ws.Range("M132").Value = -4228.750100000001
$ws.Range("N132").Value = -12270.5

$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head
$ws.Range("H2").Value = 377250
$ws.Range("I2").Value = 500000
$ws.Range("J2").Value = 9000
$ws.Range("K2").Value = 500000
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = -499888
$ws.Range("N2").Value = -9224

# Row 43: Subordinate Clause
$ws.Range("H43").Value = 20000000
$ws.Range("I43").Value = 20000000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 20000000
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("M43").Value = -19999807

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 1326.1111
$ws.Range("I68").Value = 1072
$ws.Range("J68").Value = 1834.3334
$ws.Range("K68").Value = 1072
$ws.Range("L68").Value = 1834.3334
$ws.Range("M68").Value = -323
$ws.Range("N68").Value = -3332.3334

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 1326.1111
$ws.Range("I71").Value = 1072
$ws.Range("J71").Value = 1834.3334
$ws.Range("K71").Value = 5360
$ws.Range("L71").Value = 9171.666999999999
$ws.Range("M71").Value = -1616
$ws.Range("N71").Value = -16659.667

# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 23811406
$ws.Range("I100").Value = 25642820
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 25642820
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -25642279
$ws.Range("N100").Value = -4082

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 6031.778
$ws.Range("I132").Value = 9911.958000000001
$ws.Range("J132").Value = 1597.2858
$ws.Range("K132").Value = 29735.874
$ws.Range("L132").Value = 4791.857400000001
$ws.Range("M132").Value = -27205.874
$ws.Range("N132").Value = -9851.857400000001

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 2509
$ws.Range("I81").Value = 2622.111
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 5244.222
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -4183.222
$ws.Range("N81").Value = -6122

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 2509
$ws.Range("I84").Value = 2622.111
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 26221.11
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -20917.11
$ws.Range("N84").Value = -30608

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2530.1667
$ws.Range("I132").Value = 2396.5454
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7189.6362
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4659.6362
$ws.Range("N132").Value = -17060
